# Populate the (previously empty) Sheet1 with a small contacts table.
# Mirrors the natural order a user would type the data in: header row
# first, then the data row, with a couple of early typos on the phone
# number corrected before settling on the final value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "First Name"
$ws.Range("B1").Value = "Last Name"
$ws.Range("C1").Value = "Phone Number"

# Data row
$ws.Range("A2").Value = "bob"
$ws.Range("B2").Value = "jill"

# A couple of corrections while entering the phone number before the
# final value was settled on.
$ws.Range("C2").Value = "bo"
$ws.Range("C2").Value = "jackson"
$ws.Range("C2").Value = "555-444-3333"
